$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text corrections ---
$ws.Range("G1").Value = "TALLA"
$ws.Range("H1").Value = "PESO"
$ws.Range("M1").Value = "REGISTRO FICHA?"
$ws.Range("A1").Value = "CODIGO FICHA"

# --- Column width adjustments (A and B) ---
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(2).ColumnWidth = 24.666666666666668

# --- Active cell selection ---
$ws.Range("D10").Select() | Out-Null
